$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (header "Förändrad") holds a date serial number for every data row
# (rows 2 through 536). The update moves that date forward by one day
# (serial 45202 -> 45203) for every row.
$range = $ws.Range("C2:C536")
$range.Value = 45203
